# Scheduled-runner update: refresh market-price / profit figures
# (columns H-N: currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2700.25
$ws.Range("J9").Value = 3433.6667
$ws.Range("L9").Value = 3433.6667
$ws.Range("N9").Value = -3771.6667
$ws.Range("H17").Value = 3085.0613
$ws.Range("J17").Value = 3399.2559
$ws.Range("L17").Value = 10197.7677
$ws.Range("N17").Value = -10533.7677
$ws.Range("H47").Value = 21591.334
$ws.Range("J47").Value = 27537
$ws.Range("L47").Value = 27537
$ws.Range("N47").Value = -29481
$ws.Range("H54").Value = 6817.143
$ws.Range("I54").Value = 4650
$ws.Range("K54").Value = 4650
$ws.Range("M54").Value = -4164
$ws.Range("H112").Value = 4250.25
$ws.Range("J112").Value = 4250.25
$ws.Range("L112").Value = 12750.75
$ws.Range("N112").Value = -14966.75
$ws.Range("H132").Value = 1112.3
$ws.Range("I132").Value = 1121.7037
$ws.Range("K132").Value = 3365.1111
$ws.Range("M132").Value = -835.1111000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14928131
$ws.Range("I32").Value = 15874949
$ws.Range("K32").Value = 15874949
$ws.Range("M32").Value = -15874662
$ws.Range("H45").Value = 2937.4736
$ws.Range("I45").Value = 938
$ws.Range("K45").Value = 938
$ws.Range("M45").Value = -561
$ws.Range("H61").Value = 8130.1313
$ws.Range("I61").Value = 8288.727999999999
$ws.Range("K61").Value = 8288.727999999999
$ws.Range("M61").Value = -8076.727999999999
$ws.Range("H74").Value = 13335661
$ws.Range("I74").Value = 17545470
$ws.Range("J74").Value = 4596.6665
$ws.Range("K74").Value = 17545470
$ws.Range("L74").Value = 4596.6665
$ws.Range("M74").Value = -17544596
$ws.Range("N74").Value = -6344.6665
$ws.Range("H77").Value = 13335661
$ws.Range("I77").Value = 17545470
$ws.Range("J77").Value = 4596.6665
$ws.Range("K77").Value = 87727350
$ws.Range("L77").Value = 22983.3325
$ws.Range("M77").Value = -87722982
$ws.Range("N77").Value = -31719.3325
$ws.Range("H132").Value = 1946.7843
$ws.Range("I132").Value = 1619.8223
$ws.Range("J132").Value = 4399
$ws.Range("K132").Value = 4859.4669
$ws.Range("L132").Value = 13197
$ws.Range("M132").Value = -2329.4669
$ws.Range("N132").Value = -18257
$ws.Range("H136").Value = 8130.1313
$ws.Range("I136").Value = 8288.727999999999
$ws.Range("K136").Value = 24866.184
$ws.Range("M136").Value = -22316.184

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 31999.5
$ws.Range("I63").Value = 31000
$ws.Range("J63").Value = 32999
$ws.Range("K63").Value = 31000
$ws.Range("L63").Value = 32999
$ws.Range("M63").Value = -30314
$ws.Range("N63").Value = -34371
$ws.Range("H66").Value = 31999.5
$ws.Range("I66").Value = 31000
$ws.Range("J66").Value = 32999
$ws.Range("K66").Value = 93000
$ws.Range("L66").Value = 98997
$ws.Range("M66").Value = -89568
$ws.Range("N66").Value = -105861
$ws.Range("H75").Value = 26167.875
$ws.Range("I75").Value = 13316.5
$ws.Range("J75").Value = 64722
$ws.Range("K75").Value = 13316.5
$ws.Range("L75").Value = 64722
$ws.Range("M75").Value = -12380.5
$ws.Range("N75").Value = -66594
$ws.Range("H78").Value = 26167.875
$ws.Range("I78").Value = 13316.5
$ws.Range("J78").Value = 64722
$ws.Range("K78").Value = 39949.5
$ws.Range("L78").Value = 194166
$ws.Range("M78").Value = -35269.5
$ws.Range("N78").Value = -203526
$ws.Range("H86").Value = 4530
$ws.Range("I86").Value = 3085.8
$ws.Range("J86").Value = 7418.4
$ws.Range("K86").Value = 3085.8
$ws.Range("L86").Value = 7418.4
$ws.Range("M86").Value = -1962.8
$ws.Range("N86").Value = -9664.4
$ws.Range("H89").Value = 4530
$ws.Range("I89").Value = 3085.8
$ws.Range("J89").Value = 7418.4
$ws.Range("K89").Value = 15429
$ws.Range("L89").Value = 37092
$ws.Range("M89").Value = -9813
$ws.Range("N89").Value = -48324
$ws.Range("H134").Value = 2157.64
$ws.Range("I134").Value = 944.1
$ws.Range("K134").Value = 2832.3
$ws.Range("M134").Value = -297.3000000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46743.375
$ws.Range("I31").Value = 1278.875
$ws.Range("J31").Value = 69475.625
$ws.Range("K31").Value = 1278.875
$ws.Range("L31").Value = 69475.625
$ws.Range("M31").Value = -983.875
$ws.Range("N31").Value = -70065.625
$ws.Range("H34").Value = 46743.375
$ws.Range("I34").Value = 1278.875
$ws.Range("J34").Value = 69475.625
$ws.Range("K34").Value = 1278.875
$ws.Range("L34").Value = 69475.625
$ws.Range("M34").Value = -1076.875
$ws.Range("N34").Value = -69879.625
$ws.Range("H62").Value = 9576.571
$ws.Range("I62").Value = 4548
$ws.Range("J62").Value = 11588
$ws.Range("K62").Value = 4548
$ws.Range("L62").Value = 11588
$ws.Range("M62").Value = -3924
$ws.Range("N62").Value = -12836
$ws.Range("H65").Value = 9576.571
$ws.Range("I65").Value = 4548
$ws.Range("J65").Value = 11588
$ws.Range("K65").Value = 22740
$ws.Range("L65").Value = 57940
$ws.Range("M65").Value = -19620
$ws.Range("N65").Value = -64180
$ws.Range("H69").Value = 41500
$ws.Range("I69").Value = 45333.332
$ws.Range("K69").Value = 45333.332
$ws.Range("M69").Value = -44584.332
$ws.Range("H72").Value = 41500
$ws.Range("I72").Value = 45333.332
$ws.Range("K72").Value = 135999.996
$ws.Range("M72").Value = -132255.996
$ws.Range("H105").Value = 9846.777
$ws.Range("I105").Value = 11251.667
$ws.Range("K105").Value = 11251.667
$ws.Range("M105").Value = -9504.666999999999
$ws.Range("H107").Value = 1052
$ws.Range("I107").Value = 1117.75
$ws.Range("K107").Value = 1117.75
$ws.Range("M107").Value = 802.25
$ws.Range("H134").Value = 5020.769
$ws.Range("I134").Value = 3248.353
$ws.Range("J134").Value = 8368.666999999999
$ws.Range("K134").Value = 9745.059000000001
$ws.Range("L134").Value = 25106.001
$ws.Range("M134").Value = -7210.059000000001
$ws.Range("N134").Value = -30176.001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 22.294117
$ws.Range("J12").Value = 12.571428
$ws.Range("L12").Value = 37.714284
$ws.Range("N12").Value = -383.714284
$ws.Range("H17").Value = 597.6
$ws.Range("I17").Value = 497
$ws.Range("K17").Value = 1491
$ws.Range("M17").Value = -1322
$ws.Range("H32").Value = 10885
$ws.Range("J32").Value = 12782
$ws.Range("L32").Value = 38346
$ws.Range("N32").Value = -38912
$ws.Range("H55").Value = 2383

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 190042.88
$ws.Range("I132").Value = 224728.66
$ws.Range("K132").Value = 674185.98
$ws.Range("M132").Value = -671655.98

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4014.4
$ws.Range("I132").Value = 2208.45
$ws.Range("K132").Value = 6625.349999999999
$ws.Range("M132").Value = -4095.349999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 469.1
$ws.Range("I113").Value = 232.33333
$ws.Range("K113").Value = 696.99999
$ws.Range("M113").Value = 1473.00001
$ws.Range("H126").Value = 2027.6389
$ws.Range("I126").Value = 1810.409
$ws.Range("J126").Value = 2369
$ws.Range("K126").Value = 5431.227000000001
$ws.Range("L126").Value = 7107
$ws.Range("M126").Value = -2961.227000000001
$ws.Range("N126").Value = -12047
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H132").Value = 6271.4053
$ws.Range("I132").Value = 2241.48
$ws.Range("K132").Value = 6724.440000000001
$ws.Range("M132").Value = -4194.440000000001
$ws.Range("H136").Value = 4572.6665
$ws.Range("I136").Value = 3207.5557
$ws.Range("K136").Value = 9622.667099999999
$ws.Range("M136").Value = -7072.667099999999
